$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-30 Monday" "2024-10-01 Tuesday"

Replace-Text "786÷7=" "108÷4="
Replace-Text "345÷2=" "552÷3="
Replace-Text "466÷4=" "452÷9="
Replace-Text "779÷7=" "122÷4="
Replace-Text "999÷3=" "729÷3="

Replace-Text "409÷6=" "506÷3="
Replace-Text "909÷6=" "393÷2="
Replace-Text "316÷9=" "894÷5="
Replace-Text "546÷7=" "898÷6="
Replace-Text "651÷2=" "679÷3="

Replace-Text "197÷6=" "154÷3="
Replace-Text "321÷3=" "951÷7="
Replace-Text "689÷9=" "758÷3="
Replace-Text "159÷8=" "749÷4="
Replace-Text "346÷3=" "357÷2="

Replace-Text "176÷2=" "843÷9="
Replace-Text "177÷2=" "390÷6="
Replace-Text "632÷6=" "380÷9="
Replace-Text "156÷5=" "391÷5="
Replace-Text "321÷9=" "482÷3="

Replace-Text "962÷3=" "618÷4="
Replace-Text "911÷8=" "439÷5="
Replace-Text "319÷8=" "762÷8="
Replace-Text "653÷7=" "619÷7="
Replace-Text "659÷8=" "844÷7="

Write-Output "Replacements complete"
